$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1918238993710692
$ws.Range("C2").Value = 0.5471698113207547
$ws.Range("J2").Value = 0.006289308176100629
$ws.Range("O2").Value = 0.003144654088050315
$ws.Range("P2").Value = 0.1540880503144654
$ws.Range("S2").Value = 0.09748427672955975
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.04395604395604396
$ws.Range("P3").Value = 0.7197802197802198
$ws.Range("S3").Value = 0.2142857142857143
$ws.Range("P4").Value = 0.6857142857142857
$ws.Range("S4").Value = 0.3142857142857143
$ws.Range("B6").Value = 0.08095238095238096
$ws.Range("D6").Value = 0.004761904761904762
$ws.Range("F6").Value = 0.04761904761904762
$ws.Range("J6").Value = 0.2380952380952381
$ws.Range("O6").Value = 0.01904761904761905
$ws.Range("Q6").Value = 0.1476190476190476
$ws.Range("R6").Value = 0.0761904761904762
$ws.Range("S6").Value = 0.3857142857142857
$ws.Range("B7").Value = 0.1183431952662722
$ws.Range("D7").Value = 0.03550295857988166
$ws.Range("F7").Value = 0.02958579881656805
$ws.Range("J7").Value = 0.1183431952662722
$ws.Range("O7").Value = 0.01183431952662722
$ws.Range("Q7").Value = 0.1420118343195266
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.4674556213017751
$ws.Range("B8").Value = 0.07798165137614679
$ws.Range("D8").Value = 0.009174311926605505
$ws.Range("F8").Value = 0.05045871559633028
$ws.Range("J8").Value = 0.09403669724770643
$ws.Range("O8").Value = 0.01834862385321101
$ws.Range("Q8").Value = 0.1880733944954129
$ws.Range("R8").Value = 0.07798165137614679
$ws.Range("S8").Value = 0.4839449541284404
$ws.Range("B9").Value = 0.1133333333333333
$ws.Range("D9").Value = 0.02
$ws.Range("E9").Value = 0.006666666666666667
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.12
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1733333333333333
$ws.Range("R9").Value = 0.1066666666666667
$ws.Range("S9").Value = 0.38
$ws.Range("B10").Value = 0.1281446540880503
$ws.Range("D10").Value = 0.01650943396226415
$ws.Range("F10").Value = 0.08333333333333333
$ws.Range("J10").Value = 0.08333333333333333
$ws.Range("O10").Value = 0.0110062893081761
$ws.Range("Q10").Value = 0.2075471698113208
$ws.Range("R10").Value = 0.07783018867924528
$ws.Range("S10").Value = 0.3922955974842767
$ws.Range("G11").Value = 0.1565836298932384
$ws.Range("J11").Value = 0.1067615658362989
$ws.Range("K11").Value = 0.2135231316725979
$ws.Range("L11").Value = 0.501779359430605
$ws.Range("S11").Value = 0.02135231316725979
$ws.Range("G12").Value = 0.7266666666666667
$ws.Range("J12").Value = 0.1466666666666667
$ws.Range("K12").Value = 0.02666666666666667
$ws.Range("L12").Value = 0.05333333333333334
$ws.Range("S12").Value = 0.04666666666666667
$ws.Range("G13").Value = 0.696969696969697
$ws.Range("J13").Value = 0.2424242424242424
$ws.Range("S13").Value = 0.06060606060606061
$ws.Range("F15").Value = 0.00975609756097561
$ws.Range("H15").Value = 0.1804878048780488
$ws.Range("I15").Value = 0.06829268292682927
$ws.Range("J15").Value = 0.424390243902439
$ws.Range("K15").Value = 0.04390243902439024
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("O15").Value = 0.04878048780487805
$ws.Range("S15").Value = 0.2146341463414634
$ws.Range("F16").Value = 0.009950248756218905
$ws.Range("H16").Value = 0.1890547263681592
$ws.Range("I16").Value = 0.05472636815920398
$ws.Range("J16").Value = 0.3980099502487562
$ws.Range("K16").Value = 0.1144278606965174
$ws.Range("M16").Value = 0.004975124378109453
$ws.Range("O16").Value = 0.03980099502487562
$ws.Range("S16").Value = 0.1890547263681592
$ws.Range("F17").Value = 0.009411764705882352
$ws.Range("H17").Value = 0.1623529411764706
$ws.Range("I17").Value = 0.08470588235294117
$ws.Range("J17").Value = 0.4588235294117647
$ws.Range("K17").Value = 0.08
$ws.Range("M17").Value = 0.007058823529411765
$ws.Range("O17").Value = 0.07294117647058823
$ws.Range("S17").Value = 0.1247058823529412
$ws.Range("F18").Value = 0.01129943502824859
$ws.Range("H18").Value = 0.1412429378531073
$ws.Range("I18").Value = 0.03954802259887006
$ws.Range("J18").Value = 0.5819209039548022
$ws.Range("K18").Value = 0.05084745762711865
$ws.Range("M18").Value = 0.02259887005649718
$ws.Range("O18").Value = 0.03954802259887006
$ws.Range("S18").Value = 0.1129943502824859
$ws.Range("F19").Value = 0.01415797317436662
$ws.Range("H19").Value = 0.2034277198211625
$ws.Range("I19").Value = 0.06035767511177347
$ws.Range("J19").Value = 0.3882265275707898
$ws.Range("K19").Value = 0.1028315946348733
$ws.Range("M19").Value = 0.01788375558867362
$ws.Range("N19").Value = 0.002980625931445604
$ws.Range("O19").Value = 0.07377049180327869
$ws.Range("S19").Value = 0.1363636363636364
